$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ---------------------------------------------------------------------------
# 1. Remove the "wall_thickness" sheet (its drawing + stray B29:B35 formulas
#    that pointed at never-populated rows on "data" go with it).
# ---------------------------------------------------------------------------
$wsWall = $wb.Worksheets.Item("wall_thickness")
$wsWall.Delete()

# ---------------------------------------------------------------------------
# 2. "data" sheet: append the pipe-identification rows that used to live only
#    on wall_thickness, and refresh the design-condition derived values
#    (p_avg 54.93 -> 55 bar, t_avg 601 -> 602 C ripple through V, mu, Vel,
#    Re, f, h_L, dP, total_pressure_drop).
# ---------------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("data")

$wsData.Range("B6").Value = 55
$wsData.Range("B7").Value = 602
$wsData.Range("B8").Value = 0.07155486108038436
$wsData.Range("B9").Value = 0.00003291492808074694
$wsData.Range("B10").Value = 74.95435310300522
$wsData.Range("B11").Value = 4497.261186180313
$wsData.Range("B12").Value = 23200232.39674906
$wsData.Range("B13").Value = 0.01110237907448234
$wsData.Range("B19").Value = 706.324736644795
$wsData.Range("B20").Value = 0.9680389992410972
$wsData.Range("B25").Value = 0.9680389992410972

$wsData.Range("A26").Value = "line_description"
$wsData.Range("B26").Value = "Main Pipe"
$wsData.Range("A27").Value = "nominal_size"
$wsData.Range("B27").Value = "N/A"
$wsData.Range("A28").Value = "material_and_schedule"
$wsData.Range("B28").Value = "N/A"

# ---------------------------------------------------------------------------
# 3. "pressure_drop" sheet: point the Line Description / Nominal Size /
#    Material and Schedule cells at the new "data" rows instead of being
#    blank placeholders (previously populated indirectly via wall_thickness).
# ---------------------------------------------------------------------------
$wsPD = $wb.Worksheets.Item("pressure_drop")
$wsPD.Range("D3").Formula = "=data!B26"
$wsPD.Range("D5").Formula = "=data!B27"
$wsPD.Range("D6").Formula = "=data!B28"

# Restore the sheet's active selection.
$wsPD.Activate()
$wsPD.Range("H14").Select()
